$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently holds 31 data rows (A2:B32). We need to push all of
# that data down by 11 rows (to A13:B43) to make room for 11 new rows of
# "backward extension" data at the top (A2:B12), then fill those new rows.
#
# A direct same-column Cut/Paste that overlaps its own source range behaves
# incorrectly in this runtime, so the move is staged through a scratch area
# (columns D:E) that doesn't overlap either the source or the destination.

# Step 1: move existing data out of the way (no overlap with source or target)
$ws.Range("A2:B32").Cut($ws.Range("D2"))

# Step 2: move it from the scratch area into its final shifted position
$ws.Range("D2:E32").Cut($ws.Range("A13"))

# Step 3: remove the now-empty scratch cells entirely so they don't linger
# as formatted-but-empty cells (which would otherwise inflate the sheet's
# used range / dimension).
$ws.Range("D2:E32").Delete()

# Step 4: write the new backward-extension rows (years 1983-1993 year-end
# values) into the freshly opened A2:B12 block.
$dates = @(30681,31047,31412,31777,32142,32508,32873,33238,33603,33969,34334)
$values = @(1.466797881812631,2.900424903011278,2.60323159784559,2.279090113735815,1.278816132757399,3.441722972972983,4.033476219636634,5.482086096613425,5.2529761904762,1.60115933832885,-0.9914767785701772)

for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 1).Value = $dates[$i]
    $ws.Cells.Item($row, 2).Value = $values[$i]
}

# Step 5: match the formatting that the rest of column A uses (bold, thin
# border, centered, YYYY-MM-DD HH:MM:SS date format) so the new date cells
# reuse the existing style instead of creating a new one.
$rngA = $ws.Range("A2:A12")
$rngA.NumberFormat = "YYYY-MM-DD HH:MM:SS"
$rngA.Font.Bold = $true
$rngA.HorizontalAlignment = -4108
$rngA.VerticalAlignment = -4160
$rngA.Borders.LineStyle = 1
